$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.435.32'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.639.84'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9979'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.89'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3788'
$ws.Range("E7").Value = '  -1.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.27'
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3497'
$ws.Range("E9").Value = '  -3.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08060'
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.219'
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9975'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.09'
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.317'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.259'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001201'
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").Value = '1.632.94'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.85'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06947'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.627'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9986'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.43'
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("D24").Value = '23.444.13'
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.965'
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.99'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.24'
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.171'
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.35'
$ws.Range("E30").Value = '  -1.94%  '
$ws.Range("D31").Value = '1.811.17'
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.807'
$ws.Range("E32").Value = '  -4.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.133'
$ws.Range("E33").Value = '  -5.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.20'
$ws.Range("E34").Value = '  -6.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9933'
$ws.Range("E35").Value = '  -5.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02685'
$ws.Range("E36").Value = '  -4.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08759'
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2423'
$ws.Range("E38").Value = '  -3.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.885'
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06808'
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.83'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6823'
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("E43").Value = '  -3.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.48'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9971'
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6344'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.236'
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.897'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07692'
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.11'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("E51").Value = '  +1.65%  '
